$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear column D for the listed rows (content removed, leaving an empty inline string cell)
$rowsToClear = @(3, 4, 5, 6, 8, 9, 10, 11, 12, 13, 14)
foreach ($r in $rowsToClear) {
    $ws.Cells.Item($r, 4).Value = ""
}

# Update E13 text value
$ws.Range("E13").Value = "Hesaba: Asgari 1 TL | Azami 851,5 TL"
